$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Tabelle1")

# "added bla to content ;)" -- write the two new text values into sheet1
$ws1.Range("A1").Value = "blablabla ^^"
$ws1.Range("G9").Value = "aethgälkawe fpwirag=?$*`"HTPQZ(ABGERF:<B"

# leave the selection where the author ended up after typing
[void]$ws1.Range("E14").Select()
